$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column S (02-jul) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("R1").Copy()
$ws1.Range("S1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("S1").Value = "02-jul"

$s1values = @(119.29, 102.71, 92.54000000000001, 92.08, 85.08, 93.47, 103.68, 107.8, 104.52, 92.73, 81.79000000000001, 76.39, 70.64, 50.29, 54.02, 71.56999999999999, 79.98, 87.61, 87.68000000000001, 119.03, 119.8, 126.16, 125.1, 108.51)

for ($i = 0; $i -lt $s1values.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 19).Value = $s1values[$i]
}

# --- Sheet "Gaz": add row 16 ---
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A16").NumberFormat = "@"
$ws2.Range("A16").Value = "2025-06-30"
$ws2.Range("A16").Style = "Normal"
$ws2.Range("B16").Value = 31.325

# --- Sheet "CO2": add row 16 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A16").NumberFormat = "@"
$ws3.Range("A16").Value = "2025-06-30"
$ws3.Range("A16").Style = "Normal"
$ws3.Range("B16").Value = 68
